# "park walk event pick up item"
# Insert a new "Park" event row (EV017 / Text) into the EventType(EV) sheet,
# just above the special "Every"/"OnlyScript" rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at row 18, pushing the existing row 18 ("Every", EV998, ...)
# and row 19 ("OnlyScript", EV999, ...) down to rows 19 and 20.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row with the new Park/EV017/Text event.
$ws.Range("A18").Value = "Park"
$ws.Range("B18").Value = "EV017"
$ws.Range("C18").Value = "Text"

# Match the author's final selection/cursor position on the new row.
$null = $ws.Range("C18").Select()
